$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 3236
$ws.Range("J3").Value = 3377
$ws.Range("D4").Value = 1952
$ws.Range("I4").Value = 1763
$ws.Range("J4").Value = 749
$ws.Range("J5").Value = 263
$ws.Range("J6").Value = 3984
$ws.Range("D7").Value = 28142
$ws.Range("I7").Value = 26208
$ws.Range("J7").Value = 11609

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J7").Value = 353
$ws.Range("I8").Value = 1542
$ws.Range("J8").Value = 746
$ws.Range("J11").Value = 169
$ws.Range("J12").Value = 24
$ws.Range("J15").Value = 134
$ws.Range("J19").Value = 361
$ws.Range("J20").Value = 245
$ws.Range("J21").Value = 20
$ws.Range("J23").Value = 118
$ws.Range("J25").Value = 67
$ws.Range("J29").Value = 672
$ws.Range("J30").Value = 51
$ws.Range("J31").Value = 87
$ws.Range("J33").Value = 505
$ws.Range("J36").Value = 170
$ws.Range("J37").Value = 370
$ws.Range("J42").Value = 462
$ws.Range("J46").Value = 42
$ws.Range("J47").Value = 86
$ws.Range("J48").Value = 114
$ws.Range("J52").Value = 319
$ws.Range("J54").Value = 221
$ws.Range("D63").Value = 334
$ws.Range("I63").Value = 217
$ws.Range("J63").Value = 51
$ws.Range("J65").Value = 307
$ws.Range("J67").Value = 413
$ws.Range("J73").Value = 101
$ws.Range("J76").Value = 164
$ws.Range("J78").Value = 154
$ws.Range("J79").Value = 347
$ws.Range("J83").Value = 268
$ws.Range("J84").Value = 102
$ws.Range("J85").Value = 529
$ws.Range("J86").Value = 66
$ws.Range("J88").Value = 120
$ws.Range("J89").Value = 135
$ws.Range("J90").Value = 135
$ws.Range("J91").Value = 131
$ws.Range("J96").Value = 127
$ws.Range("J99").Value = 167
$ws.Range("D101").Value = 28142
$ws.Range("I101").Value = 26208
$ws.Range("J101").Value = 11609

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 130
$ws.Range("J3").Value = 200
$ws.Range("J6").Value = 148
$ws.Range("J7").Value = 529

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 141
$ws.Range("J7").Value = 319

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 223
$ws.Range("J3").Value = 235
$ws.Range("I4").Value = 94
$ws.Range("J4").Value = 39
$ws.Range("J6").Value = 225
$ws.Range("I7").Value = 1542
$ws.Range("J7").Value = 746

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 105
$ws.Range("J7").Value = 353

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 45
$ws.Range("J4").Value = 14
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 130
$ws.Range("J5").Value = 16
$ws.Range("J7").Value = 370

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 57
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 87
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 413

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J3").Value = 23
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 102

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 91
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 79
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 162
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 505

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 205
$ws.Range("J3").Value = 231
$ws.Range("J6").Value = 169
$ws.Range("J7").Value = 672

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 104
$ws.Range("J7").Value = 361

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 96
$ws.Range("J6").Value = 229
$ws.Range("J7").Value = 462

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 51
$ws.Range("J7").Value = 154

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 9

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 61
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 101
$ws.Range("J7").Value = 347

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 76
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 245

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 42
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 101

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 27
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 31
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 24

Write-Host "Applied all changes"